# Update countries & provincias Spain
# Refresh the COVID-19 country leaderboard: updates the "last updated" timestamp,
# refreshes case counters for countries whose figures changed, and appends a new
# country row ("Guinea-Bisau") that entered the table, shifting everything below
# it down by one row (final row count goes from 201 to 202 data+header rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 title text: time of data refresh (12:46 -> 13:16)
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 13:16"

# Row -> (Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# This is the final state of every data row (rows 4-202) after the refresh.
$countryRows = @{
  4 = @("China", 81218, 47, 73650, 4287, 1399, 4, 3281)
  5 = @("Italia", 69176, 0, 8326, 54030, 3393, 0, 6820)
  6 = @("Estados Unidos", 54968, 112, 379, 53805, 1175, 4, 784)
  7 = @("España", 47610, 5552, 5367, 38809, 2636, 443, 3434)
  8 = @("Alemania", 34009, 1018, 3532, 30305, 23, 13, 172)
  9 = @("Iran", 27017, 2206, 9625, 15315, 0, 143, 2077)
  10 = @("Francia", 22304, 0, 3281, 17923, 2516, 0, 1100)
  11 = @("Suiza", 10456, 579, 131, 10180, 141, 23, 145)
  12 = @("Corea del Sur", 9137, 100, 3730, 5281, 59, 6, 126)
  13 = @("Reino Unido", 8077, 0, 135, 7520, 20, 0, 422)
  14 = @("Paises Bajos", 5560, 0, 2, 5282, 435, 0, 276)
  15 = @("Austria", 5499, 216, 9, 5460, 26, 2, 30)
  16 = @("Belgica", 4937, 668, 547, 4212, 474, 56, 178)
  17 = @("Noruega", 2902, 36, 6, 2883, 57, 1, 13)
  18 = @("Canada", 2792, 0, 112, 2654, 1, 0, 26)
  19 = @("Australia", 2431, 114, 118, 2304, 11, 1, 9)
  20 = @("Portugal", 2362, 0, 22, 2307, 48, 0, 33)
  21 = @("Suecia", 2318, 19, 16, 2261, 144, 1, 41)
  22 = @("Brasil", 2271, 24, 2, 2222, 18, 1, 47)
  23 = @("Israel", 2170, 240, 58, 2107, 37, 2, 5)
  24 = @("Turquia", 1872, 0, 0, 1828, 0, 0, 44)
  25 = @("Malasia", 1796, 172, 199, 1578, 64, 3, 19)
  26 = @("Dinamarca", 1715, 124, 1, 1680, 69, 2, 34)
  27 = @("Chequia", 1497, 103, 10, 1482, 19, 2, 5)
  28 = @("Irlanda", 1329, 0, 5, 1317, 29, 0, 7)
  29 = @("Japon", 1193, 0, 285, 865, 54, 0, 43)
  30 = @("Luxemburgo", 1099, 0, 6, 1085, 3, 0, 8)
  31 = @("Ecuador", 1082, 0, 3, 1052, 2, 0, 27)
  32 = @("Pakistan", 1000, 28, 21, 971, 5, 1, 8)
  33 = @("Tailandia", 934, 107, 70, 860, 11, 0, 4)
  34 = @("Polonia", 927, 26, 1, 914, 3, 2, 12)
  35 = @("Chile", 922, 0, 17, 903, 7, 0, 2)
  36 = @("Rumania", 906, 112, 86, 807, 18, 1, 13)
  37 = @("Finlandia", 853, 61, 10, 840, 11, 2, 3)
  38 = @("Indonesia", 790, 104, 31, 701, 0, 3, 58)
  39 = @("Arabia Saudita", 767, 0, 28, 738, 0, 0, 1)
  40 = @("Grecia", 743, 0, 29, 694, 35, 0, 20)
  41 = @("Crucero", 712, 0, 587, 115, 15, 0, 10)
  42 = @("Sudafrica", 709, 155, 12, 697, 2, 0, 0)
  43 = @("Rusia", 658, 163, 29, 628, 8, 0, 1)
  44 = @("Islandia", 648, 0, 51, 595, 13, 0, 2)
  45 = @("Filipinas", 636, 84, 26, 572, 1, 3, 38)
  46 = @("India", 562, 26, 40, 512, 0, 0, 10)
  47 = @("Singapur", 558, 0, 156, 400, 17, 0, 2)
  48 = @("Eslovenia", 528, 48, 10, 513, 14, 1, 5)
  49 = @("Catar", 526, 0, 41, 485, 6, 0, 0)
  50 = @("Panama", 443, 0, 1, 436, 33, 0, 6)
  51 = @("Barein", 419, 27, 177, 239, 2, 0, 3)
  52 = @("Croacia", 418, 36, 16, 401, 6, 0, 1)
  53 = @("Peru", 416, 0, 1, 408, 9, 0, 7)
  54 = @("Hong Kong", 410, 23, 102, 304, 4, 0, 4)
  55 = @("Mexico", 405, 38, 4, 396, 1, 1, 5)
  56 = @("Estonia", 404, 35, 8, 396, 5, 0, 0)
  57 = @("Egipto", 402, 0, 80, 302, 0, 0, 20)
  58 = @("Argentina", 387, 0, 52, 329, 0, 0, 6)
  59 = @("Colombia", 378, 0, 6, 369, 0, 0, 3)
  60 = @("Libano", 333, 15, 8, 321, 4, 0, 4)
  61 = @("Irak", 316, 0, 75, 214, 0, 0, 27)
  62 = @("Republica Dominicana", 312, 0, 3, 303, 0, 0, 6)
  63 = @("Serbia", 303, 0, 15, 284, 21, 1, 4)
  64 = @("Armenia", 265, 16, 16, 249, 6, 0, 0)
  65 = @("Argelia", 264, 0, 65, 180, 0, 0, 19)
  66 = @("Lituania", 255, 46, 1, 250, 1, 2, 4)
  67 = @("Emiratos Arabes Unidos", 248, 0, 45, 201, 2, 0, 2)
  68 = @("Taiwan", 235, 19, 29, 204, 0, 0, 2)
  69 = @("Hungria", 226, 39, 21, 195, 6, 1, 10)
  70 = @("Letonia", 221, 24, 1, 220, 0, 0, 0)
  71 = @("Bulgaria", 220, 2, 4, 213, 8, 0, 3)
  72 = @("Eslovaquia", 216, 12, 7, 209, 2, 0, 0)
  73 = @("Nueva Zelanda", 205, 0, 22, 183, 0, 0, 0)
  74 = @("Kuwait", 195, 4, 43, 152, 6, 0, 0)
  75 = @("Uruguay", 189, 0, 0, 189, 3, 0, 0)
  76 = @("San Marino", 187, 0, 4, 162, 12, 0, 21)
  77 = @("Republica de Macedonia", 177, 29, 1, 174, 1, 0, 2)
  78 = @("Costa Rica", 177, 0, 2, 173, 4, 0, 2)
  79 = @("Marruecos", 170, 0, 6, 159, 1, 0, 5)
  80 = @("Bosnia y Herzegovina", 168, 0, 2, 163, 1, 0, 3)
  81 = @("Principado de Andorra", 164, 0, 1, 162, 7, 0, 1)
  82 = @("Jordania", 154, 0, 1, 153, 0, 0, 0)
  83 = @("Albania", 146, 23, 17, 124, 3, 0, 5)
  84 = @("Vietnam", 134, 0, 17, 117, 3, 0, 0)
  85 = @("Islas Feroe", 132, 10, 38, 94, 2, 0, 0)
  86 = @("Malta", 129, 19, 2, 127, 1, 0, 0)
  87 = @("Moldavia", 125, 0, 2, 122, 20, 0, 1)
  88 = @("Republica de Chipre", 124, 0, 3, 118, 3, 0, 3)
  89 = @("Tunez", 119, 5, 1, 114, 11, 0, 4)
  90 = @("Burkina Faso", 114, 0, 7, 103, 0, 0, 4)
  91 = @("Ucrania", 113, 11, 1, 108, 0, 1, 4)
  92 = @("Brunei", 109, 5, 2, 107, 1, 0, 0)
  93 = @("Sri Lanka", 102, 0, 3, 99, 2, 0, 0)
  94 = @("Senegal", 99, 13, 9, 90, 0, 0, 0)
  95 = @("Oman", 99, 15, 17, 82, 0, 0, 0)
  96 = @("Reunion", 94, 0, 1, 93, 0, 0, 0)
  97 = @("Camboya", 93, 2, 4, 89, 1, 0, 0)
  98 = @("Venezuela", 91, 7, 15, 76, 2, 0, 0)
  99 = @("Azerbaiyan", 87, 0, 10, 76, 6, 0, 1)
  100 = @("Bielorrusia", 86, 5, 29, 57, 2, 0, 0)
  101 = @("Kazajistan", 80, 8, 0, 80, 0, 0, 0)
  102 = @("Afganistan", 79, 5, 2, 75, 0, 1, 2)
  103 = @("Guadalupe", 73, 0, 0, 72, 4, 0, 1)
  104 = @("Costa de Marfil", 73, 0, 2, 71, 0, 0, 0)
  105 = @("Georgia", 73, 3, 10, 63, 1, 0, 0)
  106 = @("Camerun", 70, 4, 2, 67, 0, 0, 1)
  107 = @("Ghana", 68, 15, 0, 66, 0, 0, 2)
  108 = @("Estado de Palestina", 62, 2, 16, 46, 0, 0, 0)
  109 = @("Trinidad yTobago", 57, 0, 0, 57, 0, 0, 0)
  110 = @("Martinica", 57, 0, 0, 56, 7, 0, 1)
  111 = @("Uzbekistan", 55, 5, 0, 55, 4, 0, 0)
  112 = @("Montenegro", 52, 5, 0, 51, 0, 0, 1)
  113 = @("Liechtenstein", 51, 0, 0, 51, 0, 0, 0)
  114 = @("Cuba", 48, 0, 1, 46, 2, 0, 1)
  115 = @("Mauricio", 48, 6, 0, 46, 1, 0, 2)
  116 = @("Nigeria", 46, 2, 2, 43, 0, 0, 1)
  117 = @("Consejo Danes para los Refugiados", 45, 0, 0, 43, 0, 0, 2)
  118 = @("Kirguistan", 44, 2, 0, 44, 0, 0, 0)
  119 = @("Ruanda", 40, 0, 0, 40, 0, 0, 0)
  120 = @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
  121 = @("Banglades", 39, 0, 7, 27, 0, 1, 5)
  122 = @("Paraguay", 37, 10, 0, 34, 1, 1, 3)
  123 = @("Mayotte", 36, 0, 0, 36, 0, 0, 0)
  124 = @("Honduras", 36, 6, 0, 36, 0, 0, 0)
  125 = @("Bolivia", 32, 3, 0, 32, 0, 0, 0)
  126 = @("Guam", 32, 0, 0, 31, 0, 0, 1)
  127 = @("Macao", 30, 1, 10, 20, 0, 0, 0)
  128 = @("Polinesia Francesa", 25, 0, 0, 25, 0, 0, 0)
  129 = @("Kenia", 25, 0, 0, 25, 0, 0, 0)
  130 = @("Jamaica", 25, 4, 2, 22, 0, 0, 1)
  131 = @("Isla de Man", 23, 0, 0, 23, 0, 0, 0)
  132 = @("Togo", 23, 3, 1, 22, 0, 0, 0)
  133 = @("Monaco", 23, 0, 1, 22, 0, 0, 0)
  134 = @("Guayana Francesa", 23, 0, 6, 17, 0, 0, 0)
  135 = @("Guatemala", 21, 0, 0, 20, 0, 0, 1)
  136 = @("Madagascar", 19, 2, 0, 19, 0, 0, 0)
  137 = @("Barbados", 18, 0, 0, 18, 0, 0, 0)
  138 = @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
  139 = @("Aruba", 17, 0, 1, 16, 0, 0, 0)
  140 = @("Gibraltar", 15, 0, 5, 10, 0, 0, 0)
  141 = @("Uganda", 14, 5, 0, 14, 0, 0, 0)
  142 = @("Nueva Caledonia", 14, 4, 0, 14, 0, 0, 0)
  143 = @("Maldivas", 13, 0, 5, 8, 0, 0, 0)
  144 = @("Etiopia", 12, 0, 0, 12, 0, 0, 0)
  145 = @("Zambia", 12, 9, 0, 12, 0, 0, 0)
  146 = @("Tanzania", 12, 0, 0, 12, 0, 0, 0)
  147 = @("Republica de Yibuti", 11, 8, 0, 11, 0, 0, 0)
  148 = @("Mongolia", 10, 0, 0, 10, 0, 0, 0)
  149 = @("Guinea Ecuatorial", 9, 0, 0, 9, 0, 0, 0)
  150 = @("El Salvador", 9, 4, 0, 9, 0, 0, 0)
  151 = @("San Martin (Parte Francesa)", 8, 0, 0, 8, 0, 0, 0)
  152 = @("Haiti", 7, 0, 0, 7, 0, 0, 0)
  153 = @("Seychelles", 7, 0, 0, 7, 0, 0, 0)
  154 = @("Dominica", 7, 0, 0, 7, 0, 0, 0)
  155 = @("Surinam", 7, 0, 0, 7, 0, 0, 0)
  156 = @("Niger", 7, 4, 0, 6, 0, 1, 1)
  157 = @("Namibia", 7, 0, 2, 5, 0, 0, 0)
  158 = @("Bermudas", 6, 0, 0, 6, 0, 0, 0)
  159 = @("Benin", 6, 0, 0, 6, 0, 0, 0)
  160 = @("Islas Caimanes", 6, 0, 0, 5, 0, 0, 1)
  161 = @("Gabon", 6, 0, 0, 5, 0, 0, 1)
  162 = @("Curazao", 6, 0, 0, 5, 0, 0, 1)
  163 = @("Fiyi", 5, 1, 0, 5, 0, 0, 0)
  164 = @("Guyana", 5, 0, 0, 4, 0, 0, 1)
  165 = @("Bahamas", 5, 0, 1, 4, 0, 0, 0)
  166 = @("Groenlandia", 5, 0, 2, 3, 0, 0, 0)
  167 = @("Suazilandia", 4, 0, 0, 4, 0, 0, 0)
  168 = @("Congo", 4, 0, 0, 4, 0, 0, 0)
  169 = @("Guinea", 4, 0, 0, 4, 0, 0, 0)
  170 = @("Santa Sede", 4, 0, 0, 4, 0, 0, 0)
  171 = @("Cabo Verde", 4, 1, 0, 3, 0, 0, 1)
  172 = @("Mozambique", 3, 0, 0, 3, 0, 0, 0)
  173 = @("Santa Lucia", 3, 0, 0, 3, 0, 0, 0)
  174 = @("Liberia", 3, 0, 0, 3, 0, 0, 0)
  175 = @("Angola", 3, 0, 0, 3, 0, 0, 0)
  176 = @("Birmania", 3, 0, 0, 3, 0, 0, 0)
  177 = @("Antigua y Barbuda", 3, 0, 0, 3, 0, 0, 0)
  178 = @("San Bartolome", 3, 0, 0, 3, 0, 0, 0)
  179 = @("Republica del Chad", 3, 0, 0, 3, 0, 0, 0)
  180 = @("Laos", 3, 1, 0, 3, 0, 0, 0)
  181 = @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0)
  182 = @("Zimbabue", 3, 0, 0, 2, 0, 0, 1)
  183 = @("Gambia", 3, 0, 0, 2, 0, 0, 1)
  184 = @("Sudan", 3, 0, 0, 2, 0, 0, 1)
  185 = @("Nepal", 3, 1, 1, 2, 0, 0, 0)
  186 = @("Nicaragua", 2, 0, 0, 2, 0, 0, 0)
  187 = @("Mauritania", 2, 0, 0, 2, 0, 0, 0)
  188 = @("San Martin (Parte Holandesa)", 2, 0, 0, 2, 0, 0, 0)
  189 = @("Guinea-Bisau", 2, 2, 0, 2, 0, 0, 0)
  190 = @("Butan", 2, 0, 0, 2, 0, 0, 0)
  191 = @("Mali", 2, 2, 0, 2, 0, 0, 0)
  192 = @("Siria", 1, 0, 0, 1, 0, 0, 0)
  193 = @("Montserrat", 1, 0, 0, 1, 0, 0, 0)
  194 = @("Eritrea", 1, 0, 0, 1, 0, 0, 0)
  195 = @("Islas Turcas y Caicos", 1, 0, 0, 1, 0, 0, 0)
  196 = @("Belice", 1, 0, 0, 1, 0, 0, 0)
  197 = @("San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0)
  198 = @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
  199 = @("Libia", 1, 0, 0, 1, 0, 0, 0)
  200 = @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0)
  201 = @("Somalia", 1, 0, 0, 1, 0, 0, 0)
  202 = @("Granada", 1, 0, 0, 0, 0, 0, 0)
}

foreach ($rowNum in $countryRows.Keys) {
  $values = $countryRows[$rowNum]
  $r = [int]$rowNum
  $ws.Cells.Item($r, 1).Value = $values[0]
  $ws.Cells.Item($r, 2).Value = $values[1]
  $ws.Cells.Item($r, 3).Value = $values[2]
  $ws.Cells.Item($r, 4).Value = $values[3]
  $ws.Cells.Item($r, 5).Value = $values[4]
  $ws.Cells.Item($r, 6).Value = $values[5]
  $ws.Cells.Item($r, 7).Value = $values[6]
  $ws.Cells.Item($r, 8).Value = $values[7]
}
